$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new value for Thursday hours on the week of 2018-03-26 (row 10)
$ws.Range("E10").Value = 4.5

# Update the active selection to match the saved state
$ws.Range("M18").Select()
